$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Customer sheet: change default/active username and move selection
# ------------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Range("A2").Value = "bartryfry@macr2.com"
$wsCustomer.Range("D5").Select()

# ------------------------------------------------------------------
# Repairer sheet: move selection to A2
# ------------------------------------------------------------------
$wsRepairer = $wb.Worksheets.Item("Repairer")
$wsRepairer.Range("A2").Select()

# ------------------------------------------------------------------
# MotorClaim_Insurer sheet: add new insurers (ICICI, HFDC ERGO,
# LIBERTY), re-sort the insurer list, extend the list/validation
# range, and add a left-border accent column (B) next to the list.
# ------------------------------------------------------------------
$wsInsurer = $wb.Worksheets.Item("MotorClaim_Insurer")

# Default dropdown selection (header row) changes to MAGMA HDI
$wsInsurer.Range("A2").Value = "MAGMA HDI"

# First add the three brand-new insurers at the bottom of the list
# (this is also how they get appended to the shared string table)
$wsInsurer.Range("A17").Value = "ICICI"
$wsInsurer.Range("A18").Value = "HFDC ERGO"
$wsInsurer.Range("A19").Value = "LIBERTY"

# Row 9 previously held the special boxed/bordered formatting (it used
# to be "NATIONAL INSURANCE HI"); that row no longer needs the accent
# once the list is re-ordered, so flatten it back to the plain style.
$wsInsurer.Range("A3").Copy()
$wsInsurer.Range("A9").PasteSpecial(-4122)

# Full alphabetised insurer list, including the three new entries
$wsInsurer.Range("A3").Value = "Bajaj Allianze General Insurance Co Ltd"
$wsInsurer.Range("A4").Value = "Chola MS General Insurance"
$wsInsurer.Range("A5").Value = "FGI MOTOR"
$wsInsurer.Range("A6").Value = "HFDC ERGO"
$wsInsurer.Range("A7").Value = "ICICI"
$wsInsurer.Range("A8").Value = "IFFCO TOKIO"
$wsInsurer.Range("A9").Value = "KOTAK MAHINDRA"
$wsInsurer.Range("A10").Value = "LIBERTY"
$wsInsurer.Range("A11").Value = "MAGMA HDI"
$wsInsurer.Range("A12").Value = "NATIONAL INSURANCE HI"
$wsInsurer.Range("A13").Value = "Royal Sundaram"
$wsInsurer.Range("A14").Value = "SBI GENERAL"
$wsInsurer.Range("A15").Value = "TATA AIG"
$wsInsurer.Range("A16").Value = "THE NEW INDIA ASSURANCE COMPANY LIMITED"
$wsInsurer.Range("A17").Value = "The Oriental Insurance Company Limited"
$wsInsurer.Range("A18").Value = "UNITED INDIA"
$wsInsurer.Range("A19").Value = "Universal Sampo General Insurance"

# Extend the blank tail of the list down to row 25 (materialise the new,
# still-empty cells with the same border formatting as the existing
# blank list cells above them).
$wsInsurer.Range("A21").Copy()
$wsInsurer.Range("A22:A25").PasteSpecial(-4122)

# Apply the existing "accent border" formatting (already used on B9/B13)
# across the whole B1:B21 column next to the list.
$template = $wsInsurer.Range("B9")
$template.Copy()
$wsInsurer.Range("B1:B21").PasteSpecial(-4122)

# Re-apply the sort over the (already ordered) list so the sheet's
# recorded sort-state grows to match the new, longer range.
$sortObj = $wsInsurer.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($wsInsurer.Range("A3:A21"))
$sortObj.SetRange($wsInsurer.Range("A3:A21"))
$sortObj.Header = -4142
$sortObj.Apply()

# Grow the header dropdown's source range to cover the longer list
$wsInsurer.Range("A2").Validation.Modify(3, 1, 1, "=`$A`$3:`$A`$25")

$wsInsurer.Activate()
$wsInsurer.Range("A2").Select()

Write-Host "done"
